$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 43 (dramaticpunch.wav), shifting
# the existing rows 43-49 down to 44-50.
$ws.Rows.Item(43).Insert()

# Fill in the new row 43 with the Twin BGM credit info.
$ws.Range("A43").Value = "lament.mp3"
$ws.Range("B43").Value = "Twin BGM"
$ws.Range("C43").Value = "4min:10sec"
$ws.Range("D43").Value = "Composer: Myuu"
$ws.Range("E43").Value = "Done"

# Give the new row its own (blank/default-looking) formatting, distinct from
# the default style used by neighboring un-styled rows.
$ws.Range("A43:E43").Interior.ColorIndex = -4142
$ws.Range("A43:E43").Borders.LineStyle = -4142

# Update the selection / scroll position on the sheet.
$ws.Range("D1").Select()
